# Update Ccl2-Ccr3 sheet with new TPM-derived values.
# Rows 2-5 get updated values; rows 6-9 (old ECs/Resolving-Mac + MuSCs/Resolving-Mac
# duplicates) are removed entirely so the sheet shrinks from A1:T9 to A1:T5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: ECs -> Resolving-Mac -------------------------------------------------
$ws.Range("D2").Value = "Resolving-Mac"
$ws.Range("G2").Value = 4.232924
$ws.Range("H2").Value = 12.698772
$ws.Range("I2").Value = 0.05792409824508498
$ws.Range("J2").Value = 0.05792409824508497
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1790523333333333
$ws.Range("N2").Value = 0.537157
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.7579149190226666
$ws.Range("R2").Value = 6.821234271204
$ws.Range("S2").Value = 0.05792409824508498
$ws.Range("T2").Value = 0.05792409824508497

# --- Row 3: FAPs / target Resolving-Mac (D3 unchanged) ---------------------------
$ws.Range("A3").Value = "FAPs"
$ws.Range("G3").Value = 23.77965533333333
$ws.Range("H3").Value = 71.338966
$ws.Range("I3").Value = 0.3254051080913003
$ws.Range("J3").Value = 0.3254051080913002
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.1790523333333333
$ws.Range("N3").Value = 0.537157
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 4.257802773295778
$ws.Range("R3").Value = 38.320224959662
$ws.Range("S3").Value = 0.3254051080913003
$ws.Range("T3").Value = 0.3254051080913002

# --- Row 4: MuSCs / target Resolving-Mac -----------------------------------------
$ws.Range("A4").Value = "MuSCs"
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("G4").Value = 7.006365333333332
$ws.Range("H4").Value = 21.019096
$ws.Range("I4").Value = 0.09587637148905993
$ws.Range("J4").Value = 0.09587637148905992
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.1790523333333333
$ws.Range("N4").Value = 0.537157
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 1.254506061119111
$ws.Range("R4").Value = 11.290554550072
$ws.Range("S4").Value = 0.09587637148905993
$ws.Range("T4").Value = 0.09587637148905992

# --- Row 5: Resolving-Mac / target Resolving-Mac ---------------------------------
$ws.Range("A5").Value = "Resolving-Mac"
$ws.Range("G5").Value = 38.058136
$ws.Range("H5").Value = 114.174408
$ws.Range("I5").Value = 0.5207944221745548
$ws.Range("J5").Value = 0.5207944221745548
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.1790523333333333
$ws.Range("N5").Value = 0.537157
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = 6.814398053117333
$ws.Range("R5").Value = 61.329582478056
$ws.Range("S5").Value = 0.5207944221745548
$ws.Range("T5").Value = 0.5207944221745548

# --- Remove old rows 6-9 entirely (sheet shrinks from A1:T9 to A1:T5) ------------
$ws.Range("A6:T9").EntireRow.Delete() | Out-Null
